$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation was recorded for Espinaca @ Vega Modelo de
# Temuco. It belongs right after the current row 304 (most recent date),
# so insert a fresh row 305 and push the existing rows 305-322 down to
# 306-323.
$ws.Rows.Item(305).Insert()

$ws.Cells.Item(305, 1).Value  = 10
$ws.Cells.Item(305, 2).Value  = 'Vega Modelo de Temuco'
$ws.Cells.Item(305, 3).Value  = 'La Araucanía'
$ws.Cells.Item(305, 4).Value  = 45265
$ws.Cells.Item(305, 5).Value  = 9
$ws.Cells.Item(305, 6).Value  = 100112012
$ws.Cells.Item(305, 7).Value  = 'Espinaca'
$ws.Cells.Item(305, 8).Value  = 'Sin especificar'
$ws.Cells.Item(305, 9).Value  = 'Primera'
$ws.Cells.Item(305, 10).Value = 30
$ws.Cells.Item(305, 11).Value = 10000
$ws.Cells.Item(305, 12).Value = 10000
$ws.Cells.Item(305, 13).Value = 10000
$ws.Cells.Item(305, 14).Value = '$/docena de paquetes'
$ws.Cells.Item(305, 15).Value = 'Región de La Araucanía'
$ws.Cells.Item(305, 16).Value = 833
$ws.Cells.Item(305, 17).Value = 12
$ws.Cells.Item(305, 18).Value = 'Hortaliza'
